$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last charge end time" (column D) for rows 2-47 to the new timestamp
$newD = 45988.344456018516
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 4).Value = $newD
}

# Refresh data rows 18-47 with updated station/terminal/time values
$ws.Cells.Item(18, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(18, 2).Value = "603号直流"
$ws.Cells.Item(18, 3).Value = 45980.250173611108
$ws.Cells.Item(19, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(19, 2).Value = "604号直流"
$ws.Cells.Item(19, 3).Value = 45985.570324074077
$ws.Cells.Item(20, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(20, 2).Value = "501号直流"
$ws.Cells.Item(20, 3).Value = 45986.210601851853
$ws.Cells.Item(21, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(21, 2).Value = "802号直流"
$ws.Cells.Item(21, 3).Value = 45986.517199074071
$ws.Cells.Item(22, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(22, 2).Value = "102号直流"
$ws.Cells.Item(22, 3).Value = 45986.537812499999
$ws.Cells.Item(23, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(23, 2).Value = "503号直流"
$ws.Cells.Item(23, 3).Value = 45986.577627314815
$ws.Cells.Item(24, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(24, 2).Value = "406号直流"
$ws.Cells.Item(24, 3).Value = 45986.586550925924
$ws.Cells.Item(25, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(25, 2).Value = "B01号直流"
$ws.Cells.Item(25, 3).Value = 45986.831574074073
$ws.Cells.Item(26, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(26, 2).Value = "A05号直流"
$ws.Cells.Item(26, 3).Value = 45987.037418981483
$ws.Cells.Item(27, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(27, 2).Value = "105号直流"
$ws.Cells.Item(27, 3).Value = 45987.098182870373
$ws.Cells.Item(28, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(28, 2).Value = "001A号直流"
$ws.Cells.Item(28, 3).Value = 45987.177488425928
$ws.Cells.Item(29, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(29, 2).Value = "805号直流"
$ws.Cells.Item(29, 3).Value = 45987.348900462966
$ws.Cells.Item(30, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(30, 2).Value = "011A号直流"
$ws.Cells.Item(30, 3).Value = 45987.412372685183
$ws.Cells.Item(31, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(31, 2).Value = "109号直流"
$ws.Cells.Item(31, 3).Value = 45987.482638888891
$ws.Cells.Item(32, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(32, 2).Value = "107号直流"
$ws.Cells.Item(32, 3).Value = 45987.524467592593
$ws.Cells.Item(33, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(33, 2).Value = "203号直流"
$ws.Cells.Item(33, 3).Value = 45987.528229166666
$ws.Cells.Item(34, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(34, 2).Value = "105号直流"
$ws.Cells.Item(34, 3).Value = 45987.533136574071
$ws.Cells.Item(35, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(35, 2).Value = "903号直流"
$ws.Cells.Item(35, 3).Value = 45987.536932870367
$ws.Cells.Item(36, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(36, 2).Value = "312号直流"
$ws.Cells.Item(36, 3).Value = 45987.53707175926
$ws.Cells.Item(37, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(37, 2).Value = "405号直流"
$ws.Cells.Item(37, 3).Value = 45987.544444444444
$ws.Cells.Item(38, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(38, 2).Value = "002A号直流"
$ws.Cells.Item(38, 3).Value = 45987.547743055555
$ws.Cells.Item(39, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(39, 2).Value = "401号直流"
$ws.Cells.Item(39, 3).Value = 45987.551886574074
$ws.Cells.Item(40, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(40, 2).Value = "101号直流"
$ws.Cells.Item(40, 3).Value = 45987.552604166667
$ws.Cells.Item(41, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(41, 2).Value = "502号直流"
$ws.Cells.Item(41, 3).Value = 45987.554328703707
$ws.Cells.Item(42, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(42, 2).Value = "111号直流"
$ws.Cells.Item(42, 3).Value = 45987.56144675926
$ws.Cells.Item(43, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(43, 2).Value = "103号直流"
$ws.Cells.Item(43, 3).Value = 45987.562800925924
$ws.Cells.Item(44, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(44, 2).Value = "009A号直流"
$ws.Cells.Item(44, 3).Value = 45987.563437500001
$ws.Cells.Item(45, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(45, 2).Value = "703号直流"
$ws.Cells.Item(45, 3).Value = 45987.563715277778
$ws.Cells.Item(46, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(46, 2).Value = "402号直流"
$ws.Cells.Item(46, 3).Value = 45987.574224537035
$ws.Cells.Item(47, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(47, 2).Value = "904号直流"
$ws.Cells.Item(47, 3).Value = 45987.576736111114

# Clear rows 48-53 (no longer have data)
for ($r = 48; $r -le 53; $r++) {
    $ws.Cells.Item($r, 1).Value = $null
    $ws.Cells.Item($r, 2).Value = $null
    $ws.Cells.Item($r, 3).Value = $null
    $ws.Cells.Item($r, 4).Value = $null
}

# Update selected cell
$ws.Range("E11").Select()
